$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 7443
$ws.Range("B2").Value = "Milena Nascimento"
$ws.Range("C2").Value = "Recursos Humanos"
$ws.Range("D2").Value = "Outros"
$ws.Range("E2").Value = 7
$ws.Range("F2").Value = 45092
$ws.Range("G2").Value = 8591.389999999999

# Row 3
$ws.Range("A3").Value = 94102
$ws.Range("B3").Value = "Ana Clara Aragão"
$ws.Range("C3").Value = "Vendas"
$ws.Range("D3").Value = "Doença"
$ws.Range("E3").Value = 8
$ws.Range("F3").Value = 45104
$ws.Range("G3").Value = 4344.86

# Row 4
$ws.Range("A4").Value = 4406
$ws.Range("B4").Value = "Srta. Maria Alice Alves"
$ws.Range("C4").Value = "Jurídico"
$ws.Range("D4").Value = "Consulta médica"
$ws.Range("E4").Value = 5
$ws.Range("F4").Value = 45090
$ws.Range("G4").Value = 3972.3

# Row 5
$ws.Range("A5").Value = 28366
$ws.Range("B5").Value = "Bruna da Costa"
$ws.Range("C5").Value = "Vendas"
$ws.Range("D5").Value = "Doença"
$ws.Range("E5").Value = 8
$ws.Range("F5").Value = 45097
$ws.Range("G5").Value = 11679.42

# Row 6
$ws.Range("A6").Value = 32474
$ws.Range("B6").Value = "João Pedro Pereira"
$ws.Range("C6").Value = "Jurídico"
$ws.Range("D6").Value = "Outros"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 45092
$ws.Range("G6").Value = 12386.31

# Row 7
$ws.Range("A7").Value = 15673
$ws.Range("B7").Value = "Bruno Barbosa"
$ws.Range("C7").Value = "Marketing"
$ws.Range("D7").Value = "Viagem de negócios"
$ws.Range("E7").Value = 4
$ws.Range("F7").Value = 45101
$ws.Range("G7").Value = 6434.46

# Row 8
$ws.Range("A8").Value = 59457
$ws.Range("B8").Value = "Luiz Henrique Barbosa"
$ws.Range("C8").Value = "Marketing"
$ws.Range("D8").Value = "Consulta médica"
$ws.Range("E8").Value = 8
$ws.Range("F8").Value = 45089
$ws.Range("G8").Value = 11533.76

# Row 9
$ws.Range("A9").Value = 39677
$ws.Range("B9").Value = "Ana Júlia da Rosa"
$ws.Range("C9").Value = "P&D"
$ws.Range("D9").Value = "Consulta médica"
$ws.Range("E9").Value = 5
$ws.Range("F9").Value = 45086
$ws.Range("G9").Value = 9494.77

# Row 10
$ws.Range("A10").Value = 62640
$ws.Range("B10").Value = "Otávio Cavalcanti"
$ws.Range("C10").Value = "Atendimento ao Cliente"
$ws.Range("D10").Value = "Consulta médica"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 45099
$ws.Range("G10").Value = 4313.52

# Row 11
$ws.Range("A11").Value = 31905
$ws.Range("B11").Value = "Theo da Cruz"
$ws.Range("C11").Value = "Engenharia"
$ws.Range("D11").Value = "Consulta médica"
$ws.Range("E11").Value = 6
$ws.Range("F11").Value = 45091
$ws.Range("G11").Value = 6453.86
